# 自动更新Excel文件 - 2025-11-28 23:12:10
# For every data row (2..99) except row 36 (malformed start-date, left untouched),
# decrement the "剩余" (remaining days, column E) by 1.
# When that countdown would reach 0, instead roll the cycle over: reset
# "剩余" back to the full cycle length ("总天", column D) and push the
# "开始时间" (start date, column F) forward by that same number of days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 99) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $total = $dCell.Value2
    $remaining = $eCell.Value2

    if ($remaining -eq $null) {
        continue
    }

    if ($remaining -le 1) {
        $eCell.Value2 = $total
        $fCell.Value2 = $fCell.Value2 + $total
    } else {
        $eCell.Value2 = $remaining - 1
    }
}

Write-Output "Updated remaining-day countdown for rows 2-$lastRow (skipped row 36)."
